$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Date serial, Error Count)
$newRows = @(
    @{ Row = 12; Date = 45971; Count = 77 },
    @{ Row = 13; Date = 45973; Count = 110 },
    @{ Row = 14; Date = 45974; Count = 65 }
)

$formatSource = $ws.Cells.Item(11, 1)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    $dateCell.Value = $r.Date
    $formatSource.Copy()
    $dateCell.PasteSpecial(-4122)  # xlPasteFormats, keep the same date number format used by the other rows

    $countCell = $ws.Cells.Item($r.Row, 2)
    $countCell.Value = $r.Count
}

# Update the selection to reflect the newly added last row, like Excel would after entry
$ws.Range("A14:B14").Select()
